$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021" column (R) by cloning the formatting already used for
# the 2020 column (Q) and then filling in the new figures.
$ws.Range("Q4:Q8").Copy()
$ws.Range("R4").PasteSpecial(-4122)

$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 47.8
$ws.Range("R6").Value = 20.7
$ws.Range("R7").Value = 9.8
$ws.Range("R8").Value = 17.3

# Match the author's resulting selection state
$ws.Range("P10").Select()
